$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the bullet "Established and led the r/CodingHelp subreddit and
#    Discord server, ..." entirely (whole paragraph, including its mark).
# ---------------------------------------------------------------------------
$target = $d.Content.Duplicate
$found = $target.Find.Execute(
    "Established and led the r/CodingHelp subreddit and Discord server, mentoring numerous individuals in JavaScript, Node.js, HTML, CSS, and SQL, fostering skill development and conducting thorough code reviews."
)
if ($found) {
    $hitPara = $target.Paragraphs.Item(1)
    foreach ($cand in $d.Paragraphs) {
        if ($cand.Range.Start -eq $hitPara.Range.Start) {
            $cand.Range.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 2) The "Developed several Discord Bots ..." bullet now talks about
#    moderating "several Discord Servers" rather than "the r/CodingHelp
#    Discord server".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " to streamline moderation tasks for the moderators of the r/CodingHelp Discord server.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " to streamline moderation tasks for the moderators of several Discord Servers.",
    2
) | Out-Null
